# BUG: Change "Chronic Renal Failure" to "Chronic Kidney Disease"
$d = $word.ActiveDocument

$find = $d.Content.Find
$find.ClearFormatting()
$find.Replacement.ClearFormatting()

# Replace "Renal Failure" with "Kidney Disease" throughout the document
# body, turning "Chronic Renal Failure" into "Chronic Kidney Disease" in
# both cells of the ICD cause-of-death table. Only the changed portion is
# replaced (the leading "Chronic " text run is left untouched), matching
# how Word's Find/Replace splits the run at the edit boundary.
$d.Content.Find.Execute("Renal Failure", $false, $false, $false,
                         $false, $false, $true, 1, $false,
                         "Kidney Disease", 2)
